$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.357.23"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "1.655.23"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "1.887.21"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "1.656.37"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "27.353.51"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("E18").Value = "  +0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.13%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("D35").Value = "1.260.21"
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.545"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.842"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  +4.93%  "
$ws.Range("E43").Value = "  +1.31%  "
$ws.Range("D44").Value = "1.797.29"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +22.79%  "
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0976"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
